# Typing practice log - add a new day's entries (rows 7 & 8) and tidy up
# the vertical alignment of the whole log table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New log rows -----------------------------------------------------
# Row 7 (Day 2): speed + problem keys
$ws.Range("H7").Value = "mainly the right hand"
$ws.Range("F7").Value = " 17-26 wpm"

# Row 8 (Day 3): day number, speed + problem keys
$ws.Range("F8").Value = "18-24 wpm"
$ws.Range("H8").Value = 'b,v.y,u,m,n,t,r,g,","'
$ws.Range("E8").Value = 3

# The new "problem keys" text for row 7 is long, shrink its font to fit
$ws.Range("H7").Font.Size = 9

# --- Re-align the whole table body to center/center -------------------
$ws.Range("E5:I22").HorizontalAlignment = -4108
$ws.Range("E5:I22").VerticalAlignment = -4108

# --- Selection cosmetic change (matches the saved workbook view) ------
$ws.Range("K6").Select()
